# "modif 2 de pagos sobre saldo"
# Update the amortization table on sheet "Hoja2": register payment #2 data,
# propagate the running balance formulas off the previous row instead of
# the fixed D7 anchor, add an interest-balance check column and a "Deuda"
# label/flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# --- Row 8: D8 now references the row above (D7) instead of the fixed $D$7 ---
$ws.Range("D8").Formula = "=D7-E8"

# --- Row 9: balance now flows from D8 (previous row) instead of fixed $D$7 ---
$ws.Range("D9").Formula = "=D8-E9"
$ws.Range("G9").Formula = "=E9+F9"

# --- Row 10: payment #2 actually registered (capital + interest amounts) ---
$ws.Range("D10").Formula = "=D9-E10"
$ws.Range("E10").Value = 5607.5
$ws.Range("F10").Value = 1392.5
$ws.Range("G10").Formula = "=E10+F10"
$ws.Range("H10").Formula = '=D9*$E$4/100/12'

# --- Row 11: interest formula now extends (shared range H9:H11) ---
$ws.Range("H11").Formula = '=D10*$E$4/100/12'

# --- New column I: "Deuda " flag/label next to row 16 and check at row 17 ---
$ws.Range("I16").Value = "Deuda "
$ws.Range("I17").Formula = "=H17-F17"

# --- New row 18: outstanding interest-debt check ---
$ws.Range("H18").Formula = "=H17-H7-H8-H9-H11"

# --- Selection moved to H12 as last edited/active cell ---
$ws.Range("H12").Select()
